$d = $word.ActiveDocument

# The document currently has a default footer only. We add a default
# header (containing the date "6/18/2013", styled with the built-in
# "Header" paragraph style) while leaving the existing footer untouched.
$section = $d.Sections.Item(1)
$headers = $section.Headers

# wdHeaderFooterPrimary = 1 -> the "default" header/footer slot.
$primaryHeader = $headers.Item(1)

# Use InsertAfter (rather than assigning Range.Text) so only the single
# "default" header part is materialized -- not the first-page/even-page
# variants as well.
$primaryHeader.Range.InsertAfter("6/18/2013")

# Apply the built-in Header paragraph style to the new paragraph, as in
# a normal Word-created header.
$headerPara = $primaryHeader.Range.Paragraphs.Item(1)
$headerPara.Style = "Header"
